$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit performs a cyclic rotation of the species-observation data across
# rows 3, 6, 7, 8, 9 (columns A,B,D,E,F,G,H,Q,R,AC), while leaving all other
# columns (C,P,S,T,U,V,W,Y,Z,AA,AB,AD,AE,AG,AT,AW,AX,AY) untouched per row.

$ws.Range("A3").Value = 111639170
$ws.Range("B3").Value = 96348
$ws.Range("D3").Value = "VU"
$ws.Range("E3").Value = 220787
$ws.Range("F3").Value = "Knärot"
$ws.Range("G3").Value = "Goodyera repens"
$ws.Range("H3").Value = "(L.) R. Br."
$ws.Range("Q3").Value = 548231.4260436196
$ws.Range("R3").Value = 6926519.619127685
$ws.Range("AC3").Value = "ca 15 plantor"

$ws.Range("A6").Value = 111639169
$ws.Range("B6").Value = 96348
$ws.Range("D6").Value = "VU"
$ws.Range("E6").Value = 220787
$ws.Range("F6").Value = "Knärot"
$ws.Range("G6").Value = "Goodyera repens"
$ws.Range("H6").Value = "(L.) R. Br."
$ws.Range("Q6").Value = 548224.5774945696
$ws.Range("R6").Value = 6926512.579557057
$ws.Range("AC6").Value = "riklig förekomst, mer än 50 plantor"

$ws.Range("A7").Value = 111639168
$ws.Range("B7").Value = 89686
$ws.Range("D7").Value = "NT"
$ws.Range("E7").Value = 658
$ws.Range("F7").Value = "Rosenticka"
$ws.Range("G7").Value = "Rhodofomes roseus"
$ws.Range("H7").Value = "(Alb. & Schwein.) Kotl. & Pouzar"
$ws.Range("Q7").Value = 548104.1391889038
$ws.Range("R7").Value = 6926477.987023209
$ws.Range("AC7").ClearContents()

$ws.Range("A8").Value = 111639174
$ws.Range("B8").Value = 96348
$ws.Range("D8").Value = "VU"
$ws.Range("E8").Value = 220787
$ws.Range("F8").Value = "Knärot"
$ws.Range("G8").Value = "Goodyera repens"
$ws.Range("H8").Value = "(L.) R. Br."
$ws.Range("Q8").Value = 547803.9854679118
$ws.Range("R8").Value = 6926147.447742103
$ws.Range("AC8").Value = "ca 6 plantor"

$ws.Range("A9").Value = 111639175
$ws.Range("B9").Value = 89686
$ws.Range("D9").Value = "NT"
$ws.Range("E9").Value = 658
$ws.Range("F9").Value = "Rosenticka"
$ws.Range("G9").Value = "Rhodofomes roseus"
$ws.Range("H9").Value = "(Alb. & Schwein.) Kotl. & Pouzar"
$ws.Range("Q9").Value = 547828.4099300706
$ws.Range("R9").Value = 6926124.660841302
$ws.Range("AC9").ClearContents()
